$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author inserted 4 new columns in front of the old "B" column, pushing the
# existing table (old B..F, plus the trailing H "notes" column) four columns to
# the right (new F..J, and L). Column A (row labels) stays put.
$ws.Range("B1:E1").EntireColumn.Insert()

# The freshly inserted columns should end up blank / unstyled, matching column A's
# width as closely as this host allows.
$ws.Range("B1:E1").Clear()
$ws.Range("B1:E1").ColumnWidth = 11.43

# New note typed into one of the newly created columns.
$ws.Range("C2").Value = "lineare SVM"

# Leave the cursor where the author ended up after the edit.
$ws.Range("C11").Select() | Out-Null

Write-Host "done"
